# Auto-generated script to update cryptos.xlsx price/volume data
# per the commit diff (Sat Apr 13 08:57:07 UTC 2024 GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.073.75'
$ws.Range("E2").Value = '  -5.09%  '
$ws.Range("D3").Value = '3.248.42'
$ws.Range("E3").Value = '  -7.72%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Formula = "'" + '594.84'
$ws.Range("E5").Value = '  -3.92%  '
$ws.Range("D6").Formula = "'" + '149.70'
$ws.Range("E6").Value = '  -13.44%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '3.240.06'
$ws.Range("E8").Value = '  -7.84%  '
$ws.Range("D9").Formula = "'" + '0.540'
$ws.Range("E9").Value = '  -11.76%  '
$ws.Range("D10").Formula = "'" + '0.170'
$ws.Range("E10").Value = '  -14.46%  '
$ws.Range("D11").Formula = "'" + '6.68'
$ws.Range("E11").Value = '  -6.08%  '
$ws.Range("D12").Formula = "'" + '0.503'
$ws.Range("E12").Value = '  -14.31%  '
$ws.Range("D13").Formula = "'" + '37.84'
$ws.Range("E13").Value = '  -18.45%  '
$ws.Range("D14").Formula = "'" + '0.0000241'
$ws.Range("E14").Value = '  -12.74%  '
$ws.Range("D15").Value = '3.770.19'
$ws.Range("E15").Value = '  -7.81%  '
$ws.Range("D16").Value = '67.098.18'
$ws.Range("E16").Value = '  -5.12%  '
$ws.Range("D17").Value = '3.254.81'
$ws.Range("E17").Value = '  -7.55%  '
$ws.Range("D18").Formula = "'" + '537.63'
$ws.Range("E18").Value = '  -11.77%  '
$ws.Range("E19").Value = '  -6.26%  '
$ws.Range("D20").Formula = "'" + '7.15'
$ws.Range("E20").Value = '  -14.43%  '
$ws.Range("D21").Formula = "'" + '15.04'
$ws.Range("E21").Value = '  -15.21%  '
$ws.Range("D22").Formula = "'" + '0.757'
$ws.Range("E22").Value = '  -14.29%  '
$ws.Range("D23").Formula = "'" + '7.83'
$ws.Range("E23").Value = '  -14.78%  '
$ws.Range("D24").Formula = "'" + '85.06'
$ws.Range("E24").Value = '  -12.98%  '
$ws.Range("D25").Formula = "'" + '13.39'
$ws.Range("E25").Value = '  -14.27%  '
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("D27").Formula = "'" + '3.24'
$ws.Range("E27").Value = '  -13.01%  '
$ws.Range("D28").Formula = "'" + '29.22'
$ws.Range("E28").Value = '  -13.19%  '
$ws.Range("D29").Formula = "'" + '7.94'
$ws.Range("E29").Value = '  -12.57%  '
$ws.Range("D30").Formula = "'" + '2.11'
$ws.Range("E30").Value = '  -17.87%  '
$ws.Range("D31").Formula = "'" + '2.62'
$ws.Range("E31").Value = '  -12.43%  '
$ws.Range("D32").Formula = "'" + '1.13'
$ws.Range("E32").Value = '  -13.05%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Formula = "'" + '6.61'
$ws.Range("E33").Value = '  -18.19%  '
$ws.Range("B34").Value = 'Bittensor'
$ws.Range("C34").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D34").Formula = "'" + '537.05'
$ws.Range("E34").Value = '  -16.25%  '
$ws.Range("D35").Formula = "'" + '5.65'
$ws.Range("E35").Value = '  -17.07%  '
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("B37").Value = 'OKB'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D37").Formula = "'" + '53.05'
$ws.Range("E37").Value = '  -6.65%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Formula = "'" + '0.0437'
$ws.Range("E38").Value = '  -9.92%  '
$ws.Range("D39").Formula = "'" + '0.0848'
$ws.Range("E39").Value = '  -14.90%  '
$ws.Range("D40").Formula = "'" + '9.08'
$ws.Range("E40").Value = '  -15.89%  '
$ws.Range("D41").Formula = "'" + '0.127'
$ws.Range("E41").Value = '  -10.62%  '
$ws.Range("D42").Value = '2.910.62'
$ws.Range("E42").Value = '  -13.15%  '
$ws.Range("D43").Formula = "'" + '2.64'
$ws.Range("E43").Value = '  -22.63%  '
$ws.Range("D44").Formula = "'" + '0.260'
$ws.Range("E44").Value = '  -16.77%  '
$ws.Range("D45").Value = '0.0₃0577'
$ws.Range("E45").Value = '  -19.71%  '
$ws.Range("D46").Formula = "'" + '2.15'
$ws.Range("E46").Value = '  -15.11%  '
$ws.Range("D47").Formula = "'" + '26.29'
$ws.Range("E47").Value = '  -17.35%  '
$ws.Range("E48").Value = '  -0.09%  '
$ws.Range("D49").Formula = "'" + '126.53'
$ws.Range("E49").Value = '  -5.65%  '
$ws.Range("D50").Formula = "'" + '2.31'
$ws.Range("E50").Value = '  -22.21%  '
$ws.Range("E51").Value = '  -13.23%  '
